# Update crypto price/volume/hour snapshot columns (symbol-list refresh).
# Values in columns D (Price) and G (Hora) are numeric-looking but are
# stored as text in the source sheet, so a leading apostrophe forces Excel
# to keep them as text instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''241.26'
$ws.Range("G2").Value = '''15'
$ws.Range("D3").Value = '''21.86'
$ws.Range("G3").Value = '''15'
$ws.Range("D4").Value = '''5.370'
$ws.Range("G4").Value = '''15'
$ws.Range("D5").Value = '''0.05704'
$ws.Range("G5").Value = '''15'
$ws.Range("D6").Value = '''3.431'
$ws.Range("G6").Value = '''15'
$ws.Range("D7").Value = '''6.292'
$ws.Range("G7").Value = '''15'
$ws.Range("D8").Value = '''0.8062'
$ws.Range("G8").Value = '''15'
$ws.Range("D9").Value = '''0.8534'
$ws.Range("G9").Value = '''15'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1434'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("G10").Value = '''15'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.07288'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("G11").Value = '''15'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = '''0.03069'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G12").Value = '''15'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03144'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("G13").Value = '''15'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09357'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("G14").Value = '''15'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = '''3.937'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("G15").Value = '''15'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '''0.001586'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("G16").Value = '''15'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = '''0.04830'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("G17").Value = '''15'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = '''0.0005856'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("G18").Value = '''15'
$ws.Range("D19").Value = '''0.006340'
$ws.Range("G19").Value = '''15'
$ws.Range("D20").Value = '''0.0009994'
$ws.Range("G20").Value = '''15'
$ws.Range("D21").Value = '''0.004065'
$ws.Range("G21").Value = '''15'
$ws.Range("D22").Value = '''0.0001504'
$ws.Range("G22").Value = '''15'
$ws.Range("D23").Value = '''3.724'
$ws.Range("G23").Value = '''15'
$ws.Range("G24").Value = '''15'
$ws.Range("D25").Value = '''0.3267'
$ws.Range("G25").Value = '''15'
$ws.Range("D26").Value = '''0.1244'
$ws.Range("G26").Value = '''15'
$ws.Range("D27").Value = '''0.0004011'
$ws.Range("G27").Value = '''15'
$ws.Range("G28").Value = '''15'
$ws.Range("G29").Value = '''15'
$ws.Range("G30").Value = '''15'
$ws.Range("G31").Value = '''15'
$ws.Range("G32").Value = '''15'
$ws.Range("G33").Value = '''15'
$ws.Range("G34").Value = '''15'
$ws.Range("G35").Value = '''15'
$ws.Range("G36").Value = '''15'
$ws.Range("G37").Value = '''15'
$ws.Range("G38").Value = '''15'
$ws.Range("G39").Value = '''15'
$ws.Range("D40").Value = '''0.03833'
$ws.Range("G40").Value = '''15'
$ws.Range("D41").Value = '''0.006772'
$ws.Range("E41").Value = '40KickTokenKICKBestin24h'
$ws.Range("G41").Value = '''15'
$ws.Range("D42").Value = '''0.1051'
$ws.Range("G42").Value = '''15'
$ws.Range("D43").Value = '''0.002807'
$ws.Range("G43").Value = '''15'
$ws.Range("D44").Value = '''0.007340'
$ws.Range("G44").Value = '''15'
$ws.Range("D45").Value = '''0.00005612'
$ws.Range("G45").Value = '''15'
$ws.Range("G46").Value = '''15'
$ws.Range("D47").Value = '''0.5816'
$ws.Range("G47").Value = '''15'
$ws.Range("D48").Value = '''0.1430'
$ws.Range("G48").Value = '''15'
$ws.Range("D49").Value = '''0.00002106'
$ws.Range("G49").Value = '''15'
$ws.Range("G50").Value = '''15'
$ws.Range("G51").Value = '''15'
